$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the passive-stiffness data values (Lichtwark-deleted values tweak)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 384.50110722503314
$ws.Range("C2").Value = 458.24148137044847
$ws.Range("D2").Value = 381.29154584488936
$ws.Range("E2").Value = 457.56741372270648

$ws.Range("B3").Value = 387.1899712143674
$ws.Range("C3").Value = 470.35786732153366
$ws.Range("D3").Value = 387.66631878104988
$ws.Range("E3").Value = 461.13483003093074

# Update selection to match the new data extent highlighted in the workbook
$ws.Range("B1:E3").Select()
